$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44469
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 29000
$ws.Range("P2").Value = 28500
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 2850
$ws.Range("T2").Value = 10

$ws.Range("D3").Value = 44484
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 25500
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 2550
$ws.Range("T3").Value = 10

$ws.Range("D4").Value = 44446
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 3200
$ws.Range("O4").Value = 3300
$ws.Range("P4").Value = 3250
$ws.Range("Q4").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R4").Value = 'Provincia del Elquí'
$ws.Range("S4").Value = 3250
$ws.Range("T4").Value = 1

$ws.Range("D5").Value = 44530
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 2000
$ws.Range("O5").Value = 2100
$ws.Range("P5").Value = 2050
$ws.Range("Q5").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 2050
$ws.Range("T5").Value = 1

$ws.Range("D6").Value = 44516
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 1900
$ws.Range("O6").Value = 2000
$ws.Range("P6").Value = 1950
$ws.Range("Q6").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 1950
$ws.Range("T6").Value = 1

$ws.Range("D7").Value = 44516
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 1700
$ws.Range("O7").Value = 1700
$ws.Range("P7").Value = 1700
$ws.Range("Q7").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R7").Value = 'Provincia de Limarí'
$ws.Range("S7").Value = 1700
$ws.Range("T7").Value = 1

$ws.Range("D8").Value = 44475
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 29000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 29500
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 2950
$ws.Range("T8").Value = 10

$ws.Range("D9").Value = 44495
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 26000
$ws.Range("O9").Value = 27000
$ws.Range("P9").Value = 26500
$ws.Range("Q9").Value = '$/bandeja 10 kilos'
$ws.Range("R9").Value = 'Provincia de Limarí'
$ws.Range("S9").Value = 2650
$ws.Range("T9").Value = 10

$ws.Range("D10").Value = 44505
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 2200
$ws.Range("O10").Value = 2200
$ws.Range("P10").Value = 2200
$ws.Range("Q10").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R10").Value = 'Provincia de Limarí'
$ws.Range("S10").Value = 2200
$ws.Range("T10").Value = 1

$ws.Range("D11").Value = 44505
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 1800
$ws.Range("O11").Value = 1800
$ws.Range("P11").Value = 1800
$ws.Range("Q11").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 1800
$ws.Range("T11").Value = 1

$ws.Range("D12").Value = 44160
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 17500
$ws.Range("Q12").Value = '$/bandeja 8 kilos'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 2188
$ws.Range("T12").Value = 8

$ws.Range("D13").Value = 44160
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("Q13").Value = '$/bandeja 8 kilos'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 1875
$ws.Range("T13").Value = 8

$ws.Range("D14").Value = 44467
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 2700
$ws.Range("O14").Value = 2800
$ws.Range("P14").Value = 2750
$ws.Range("Q14").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 2750
$ws.Range("T14").Value = 1

$ws.Range("D15").Value = 44467
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 2500
$ws.Range("O15").Value = 2500
$ws.Range("P15").Value = 2500
$ws.Range("Q15").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R15").Value = 'Provincia de Limarí'
$ws.Range("S15").Value = 2500
$ws.Range("T15").Value = 1

$ws.Range("D16").Value = 44517
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 27000
$ws.Range("P16").Value = 26000
$ws.Range("Q16").Value = '$/bandeja 10 kilos'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 2600
$ws.Range("T16").Value = 10

$ws.Range("D17").Value = 44483
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 2600
$ws.Range("O17").Value = 2600
$ws.Range("P17").Value = 2600
$ws.Range("Q17").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 2600
$ws.Range("T17").Value = 1

$ws.Range("D18").Value = 44483
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 2400
$ws.Range("O18").Value = 2400
$ws.Range("P18").Value = 2400
$ws.Range("Q18").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R18").Value = 'Provincia de Limarí'
$ws.Range("S18").Value = 2400
$ws.Range("T18").Value = 1

$ws.Range("D19").Value = 44488
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 25000
$ws.Range("O19").Value = 26000
$ws.Range("P19").Value = 25600
$ws.Range("Q19").Value = '$/bandeja 10 kilos'
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 2560
$ws.Range("T19").Value = 10

$ws.Range("D20").Value = 44461
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 29000
$ws.Range("O20").Value = 30000
$ws.Range("P20").Value = 29500
$ws.Range("Q20").Value = '$/bandeja 10 kilos'
$ws.Range("R20").Value = 'Provincia de Limarí'
$ws.Range("S20").Value = 2950
$ws.Range("T20").Value = 10

$ws.Range("D21").Value = 44454
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 30000
$ws.Range("O21").Value = 31000
$ws.Range("P21").Value = 30500
$ws.Range("Q21").Value = '$/bandeja 10 kilos'
$ws.Range("R21").Value = 'Provincia de Limarí'
$ws.Range("S21").Value = 3050
$ws.Range("T21").Value = 10

$ws.Range("D22").Value = 44491
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = 25000
$ws.Range("O22").Value = 26000
$ws.Range("P22").Value = 25467
$ws.Range("Q22").Value = '$/bandeja 10 kilos'
$ws.Range("R22").Value = 'Provincia de Limarí'
$ws.Range("S22").Value = 2547
$ws.Range("T22").Value = 10

$ws.Range("D23").Value = 44462
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 2900
$ws.Range("O23").Value = 3000
$ws.Range("P23").Value = 2950
$ws.Range("Q23").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R23").Value = 'Provincia de Limarí'
$ws.Range("S23").Value = 2950
$ws.Range("T23").Value = 1

$ws.Range("D24").Value = 44462
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 2600
$ws.Range("O24").Value = 2600
$ws.Range("P24").Value = 2600
$ws.Range("Q24").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R24").Value = 'Provincia de Limarí'
$ws.Range("S24").Value = 2600
$ws.Range("T24").Value = 1

$ws.Range("D25").Value = 44511
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 80
$ws.Range("N25").Value = 25000
$ws.Range("O25").Value = 26000
$ws.Range("P25").Value = 25375
$ws.Range("Q25").Value = '$/bandeja 10 kilos'
$ws.Range("R25").Value = 'Provincia de Limarí'
$ws.Range("S25").Value = 2538
$ws.Range("T25").Value = 10
